$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 225.875  # H2: 199.45454 -> 225.875
$ws.Cells.Item(2, 9).Value = 137  # I2: 124.28571 -> 137.0
$ws.Cells.Item(2, 10).Value = 374  # J2: 331.0 -> 374.0
$ws.Cells.Item(2, 11).Value = 137  # K2: 124.28571 -> 137.0
$ws.Cells.Item(2, 12).Value = 374  # L2: 331.0 -> 374.0
$ws.Cells.Item(2, 13).Value = -24  # M2: -11.28570999999999 -> -24.0
$ws.Cells.Item(2, 14).Value = -600  # N2: -557.0 -> -600.0
$ws.Cells.Item(19, 8).Value = 500  # H19: 674.25 -> 500.0
$ws.Cells.Item(19, 9).Value = 500  # I19: 499.0 -> 500.0
$ws.Cells.Item(19, 10).Value = 0  # J19: 1200.0 -> 0.0
$ws.Cells.Item(19, 11).Value = 500  # K19: 499.0 -> 500.0
$ws.Cells.Item(19, 12).Value = 0  # L19: 1200.0 -> 0.0
$ws.Cells.Item(19, 13).Value = -325  # M19: -324.0 -> -325.0
$ws.Cells.Item(19, 14).ClearContents()  # N19: -1550.0 -> (removed)
$ws.Cells.Item(40, 8).Value = 7504  # H40: 7789.0713 -> 7504.0
$ws.Cells.Item(40, 10).Value = 9899.2  # J40: 9785.0 -> 9899.2
$ws.Cells.Item(40, 12).Value = 9899.2  # L40: 9785.0 -> 9899.2
$ws.Cells.Item(40, 14).Value = -10249.2  # N40: -10135.0 -> -10249.2
$ws.Cells.Item(43, 8).Value = 2476.6  # H43: 2043.5 -> 2476.6
$ws.Cells.Item(43, 9).Value = 2742.5  # I43: 1840.5 -> 2742.5
$ws.Cells.Item(43, 10).Value = 2299.3333  # J43: 2449.5 -> 2299.3333
$ws.Cells.Item(43, 11).Value = 2742.5  # K43: 1840.5 -> 2742.5
$ws.Cells.Item(43, 12).Value = 2299.3333  # L43: 2449.5 -> 2299.3333
$ws.Cells.Item(43, 13).Value = -2673.5  # M43: -1771.5 -> -2673.5
$ws.Cells.Item(43, 14).Value = -2437.3333  # N43: -2587.5 -> -2437.3333
$ws.Cells.Item(48, 8).Value = 0  # H48: 1500.0 -> 0.0
$ws.Cells.Item(48, 10).Value = 0  # J48: 1500.0 -> 0.0
$ws.Cells.Item(48, 12).Value = 0  # L48: 4500.0 -> 0.0
$ws.Cells.Item(48, 14).ClearContents()  # N48: -5084.0 -> (removed)
$ws.Cells.Item(56, 8).Value = 0  # H56: 1500.0 -> 0.0
$ws.Cells.Item(56, 10).Value = 0  # J56: 1500.0 -> 0.0
$ws.Cells.Item(56, 12).Value = 0  # L56: 4500.0 -> 0.0
$ws.Cells.Item(56, 14).ClearContents()  # N56: -5568.0 -> (removed)
$ws.Cells.Item(69, 8).Value = 7580  # H69: 7550.1665 -> 7580.0
$ws.Cells.Item(69, 10).Value = 7580  # J69: 7550.1665 -> 7580.0
$ws.Cells.Item(69, 12).Value = 22740  # L69: 22650.4995 -> 22740.0
$ws.Cells.Item(69, 14).Value = -24488  # N69: -24398.4995 -> -24488.0
$ws.Cells.Item(72, 8).Value = 7580  # H72: 7550.1665 -> 7580.0
$ws.Cells.Item(72, 10).Value = 7580  # J72: 7550.1665 -> 7580.0
$ws.Cells.Item(72, 12).Value = 68220  # L72: 67951.4985 -> 68220.0
$ws.Cells.Item(72, 14).Value = -76956  # N72: -76687.4985 -> -76956.0
$ws.Cells.Item(80, 8).Value = 187.5  # H80: 0.0 -> 187.5
$ws.Cells.Item(80, 9).Value = 116.666664  # I80: 0.0 -> 116.666664
$ws.Cells.Item(80, 10).Value = 400  # J80: 0.0 -> 400.0
$ws.Cells.Item(80, 11).Value = 349.999992  # K80: 0.0 -> 349.999992
$ws.Cells.Item(80, 12).Value = 1200  # L80: 0.0 -> 1200.0
$ws.Cells.Item(80, 13).Value = 648.000008  # M80: None -> 648.000008
$ws.Cells.Item(80, 14).Value = -3196  # N80: None -> -3196.0
$ws.Cells.Item(83, 8).Value = 187.5  # H83: 0.0 -> 187.5
$ws.Cells.Item(83, 9).Value = 116.666664  # I83: 0.0 -> 116.666664
$ws.Cells.Item(83, 10).Value = 400  # J83: 0.0 -> 400.0
$ws.Cells.Item(83, 11).Value = 1049.999976  # K83: 0.0 -> 1049.999976
$ws.Cells.Item(83, 12).Value = 3600  # L83: 0.0 -> 3600.0
$ws.Cells.Item(83, 13).Value = 3942.000024  # M83: None -> 3942.000024
$ws.Cells.Item(83, 14).Value = -13584  # N83: None -> -13584.0
$ws.Cells.Item(98, 8).Value = 898.5  # H98: 1318.8 -> 898.5
$ws.Cells.Item(98, 10).Value = 1199.5  # J98: 1799.6666 -> 1199.5
$ws.Cells.Item(98, 12).Value = 1199.5  # L98: 1799.6666 -> 1199.5
$ws.Cells.Item(98, 14).Value = -4195.5  # N98: -4795.6666 -> -4195.5
$ws.Cells.Item(100, 8).Value = 1120.0834  # H100: 1329.1 -> 1120.0834
$ws.Cells.Item(100, 9).Value = 1068.25  # I100: 1257.4286 -> 1068.25
$ws.Cells.Item(100, 10).Value = 1223.75  # J100: 1496.3334 -> 1223.75
$ws.Cells.Item(100, 11).Value = 1068.25  # K100: 1257.4286 -> 1068.25
$ws.Cells.Item(100, 12).Value = 1223.75  # L100: 1496.3334 -> 1223.75
$ws.Cells.Item(100, 13).Value = -527.25  # M100: -716.4286 -> -527.25
$ws.Cells.Item(100, 14).Value = -2305.75  # N100: -2578.3334 -> -2305.75
$ws.Cells.Item(122, 8).Value = 898.5  # H122: 1318.8 -> 898.5
$ws.Cells.Item(122, 10).Value = 1199.5  # J122: 1799.6666 -> 1199.5
$ws.Cells.Item(122, 12).Value = 3598.5  # L122: 5398.9998 -> 3598.5
$ws.Cells.Item(122, 14).Value = -8498.5  # N122: -10298.9998 -> -8498.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 445.4  # H4: 457.6 -> 445.4
$ws.Cells.Item(4, 9).Value = 277  # I4: 287.5 -> 277.0
$ws.Cells.Item(4, 10).Value = 557.6667  # J4: 571.0 -> 557.6667
$ws.Cells.Item(4, 11).Value = 277  # K4: 287.5 -> 277.0
$ws.Cells.Item(4, 12).Value = 557.6667  # L4: 571.0 -> 557.6667
$ws.Cells.Item(4, 13).Value = -161  # M4: -171.5 -> -161.0
$ws.Cells.Item(4, 14).Value = -789.6667  # N4: -803.0 -> -789.6667
$ws.Cells.Item(32, 8).Value = 1039.9615  # H32: 1019.4815 -> 1039.9615
$ws.Cells.Item(32, 9).Value = 982.13043  # I32: 961.5 -> 982.13043
$ws.Cells.Item(32, 11).Value = 982.13043  # K32: 961.5 -> 982.13043
$ws.Cells.Item(32, 13).Value = -695.13043  # M32: -674.5 -> -695.13043
$ws.Cells.Item(63, 8).Value = 6522  # H63: 5934.8335 -> 6522.0
$ws.Cells.Item(63, 9).Value = 2902.5  # I63: 2934.6667 -> 2902.5
$ws.Cells.Item(63, 11).Value = 2902.5  # K63: 2934.6667 -> 2902.5
$ws.Cells.Item(63, 13).Value = -2216.5  # M63: -2248.6667 -> -2216.5
$ws.Cells.Item(66, 8).Value = 6522  # H66: 5934.8335 -> 6522.0
$ws.Cells.Item(66, 9).Value = 2902.5  # I66: 2934.6667 -> 2902.5
$ws.Cells.Item(66, 11).Value = 14512.5  # K66: 14673.3335 -> 14512.5
$ws.Cells.Item(66, 13).Value = -11080.5  # M66: -11241.3335 -> -11080.5
$ws.Cells.Item(97, 8).Value = 1153.4736  # H97: 1118.3 -> 1153.4736
$ws.Cells.Item(97, 9).Value = 1151.375  # I97: 1110.1177 -> 1151.375
$ws.Cells.Item(97, 11).Value = 1151.375  # K97: 1110.1177 -> 1151.375
$ws.Cells.Item(97, 13).Value = -655.375  # M97: -614.1177 -> -655.375

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4644.273  # H86: 4053.2727 -> 4644.273
$ws.Cells.Item(86, 10).Value = 7165.5  # J86: 6082.0 -> 7165.5
$ws.Cells.Item(86, 12).Value = 7165.5  # L86: 6082.0 -> 7165.5
$ws.Cells.Item(86, 14).Value = -9411.5  # N86: -8328.0 -> -9411.5
$ws.Cells.Item(89, 8).Value = 4644.273  # H89: 4053.2727 -> 4644.273
$ws.Cells.Item(89, 10).Value = 7165.5  # J89: 6082.0 -> 7165.5
$ws.Cells.Item(89, 12).Value = 35827.5  # L89: 30410.0 -> 35827.5
$ws.Cells.Item(89, 14).Value = -47059.5  # N89: -41642.0 -> -47059.5
$ws.Cells.Item(94, 8).Value = 982.8461  # H94: 1037.4615 -> 982.8461
$ws.Cells.Item(94, 9).Value = 743.4545  # I94: 808.0 -> 743.4545
$ws.Cells.Item(94, 11).Value = 743.4545  # K94: 808.0 -> 743.4545
$ws.Cells.Item(94, 13).Value = -292.4545000000001  # M94: -357.0 -> -292.4545000000001
$ws.Cells.Item(107, 8).Value = 6044.5  # H107: 5768.5454 -> 6044.5
$ws.Cells.Item(107, 9).Value = 5828.3335  # I107: 5546.4 -> 5828.3335
$ws.Cells.Item(107, 11).Value = 5828.3335  # K107: 5546.4 -> 5828.3335
$ws.Cells.Item(107, 13).Value = -3908.3335  # M107: -3626.4 -> -3908.3335
$ws.Cells.Item(140, 8).Value = 78195  # H140: 88333.336 -> 78195.0
$ws.Cells.Item(140, 10).Value = 78195  # J140: 88333.336 -> 78195.0
$ws.Cells.Item(140, 12).Value = 78195  # L140: 88333.336 -> 78195.0
$ws.Cells.Item(140, 14).Value = -88555  # N140: -98693.336 -> -88555.0

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1700.8667  # H22: 1707.4667 -> 1700.8667
$ws.Cells.Item(22, 9).Value = 1278.9231  # I22: 1286.5385 -> 1278.9231
$ws.Cells.Item(22, 11).Value = 1278.9231  # K22: 1286.5385 -> 1278.9231
$ws.Cells.Item(22, 13).Value = -928.9231  # M22: -936.5385000000001 -> -928.9231
$ws.Cells.Item(62, 8).Value = 2973.8  # H62: 2420.5 -> 2973.8
$ws.Cells.Item(62, 9).Value = 2904.75  # I62: 2302.0 -> 2904.75
$ws.Cells.Item(62, 11).Value = 2904.75  # K62: 2302.0 -> 2904.75
$ws.Cells.Item(62, 13).Value = -2280.75  # M62: -1678.0 -> -2280.75
$ws.Cells.Item(65, 8).Value = 2973.8  # H65: 2420.5 -> 2973.8
$ws.Cells.Item(65, 9).Value = 2904.75  # I65: 2302.0 -> 2904.75
$ws.Cells.Item(65, 11).Value = 14523.75  # K65: 11510.0 -> 14523.75
$ws.Cells.Item(65, 13).Value = -11403.75  # M65: -8390.0 -> -11403.75
$ws.Cells.Item(107, 8).Value = 955.94116  # H107: 1013.86664 -> 955.94116
$ws.Cells.Item(107, 9).Value = 384.81818  # I107: 367.91666 -> 384.81818
$ws.Cells.Item(107, 10).Value = 2003  # J107: 3597.6667 -> 2003.0
$ws.Cells.Item(107, 11).Value = 384.81818  # K107: 367.91666 -> 384.81818
$ws.Cells.Item(107, 12).Value = 2003  # L107: 3597.6667 -> 2003.0
$ws.Cells.Item(107, 13).Value = 1535.18182  # M107: 1552.08334 -> 1535.18182
$ws.Cells.Item(107, 14).Value = -5843  # N107: -7437.6667 -> -5843.0
$ws.Cells.Item(132, 8).Value = 3326  # H132: 2003.375 -> 3326.0
$ws.Cells.Item(132, 9).Value = 3326  # I132: 2003.375 -> 3326.0
$ws.Cells.Item(132, 11).Value = 9978  # K132: 6010.125 -> 9978.0
$ws.Cells.Item(132, 13).Value = -7448  # M132: -3480.125 -> -7448.0

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 6710.8184  # H55: 6782.2856 -> 6710.8184
$ws.Cells.Item(55, 9).Value = 704.5  # I55: 566.3333 -> 704.5
$ws.Cells.Item(55, 10).Value = 7311.45  # J55: 7818.278 -> 7311.45
$ws.Cells.Item(55, 11).Value = 2113.5  # K55: 1698.9999 -> 2113.5
$ws.Cells.Item(55, 12).Value = 21934.35  # L55: 23454.834 -> 21934.35
$ws.Cells.Item(55, 13).Value = -1936.5  # M55: -1521.9999 -> -1936.5
$ws.Cells.Item(55, 14).Value = -22288.35  # N55: -23808.834 -> -22288.35
$ws.Cells.Item(64, 8).Value = 3500  # H64: 0.0 -> 3500.0
$ws.Cells.Item(64, 9).Value = 3500  # I64: 0.0 -> 3500.0
$ws.Cells.Item(64, 11).Value = 10500  # K64: 0.0 -> 10500.0
$ws.Cells.Item(64, 13).Value = -10230  # M64: None -> -10230.0
$ws.Cells.Item(67, 8).Value = 3500  # H67: 0.0 -> 3500.0
$ws.Cells.Item(67, 9).Value = 3500  # I67: 0.0 -> 3500.0
$ws.Cells.Item(67, 11).Value = 10500  # K67: 0.0 -> 10500.0
$ws.Cells.Item(67, 13).Value = -9564  # M67: None -> -9564.0
$ws.Cells.Item(102, 8).Value = 0  # H102: 1750.0 -> 0.0
$ws.Cells.Item(102, 10).Value = 0  # J102: 1750.0 -> 0.0
$ws.Cells.Item(102, 12).Value = 0  # L102: 5250.0 -> 0.0
$ws.Cells.Item(102, 14).ClearContents()  # N102: -10118.0 -> (removed)

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 13560.667  # H57: 14597.667 -> 13560.667
$ws.Cells.Item(57, 10).Value = 13560.667  # J57: 14597.667 -> 13560.667
$ws.Cells.Item(57, 12).Value = 13560.667  # L57: 14597.667 -> 13560.667
$ws.Cells.Item(57, 14).Value = -15200.667  # N57: -16237.667 -> -15200.667
$ws.Cells.Item(70, 8).Value = 8000  # H70: 5999.5 -> 8000.0
$ws.Cells.Item(70, 9).Value = 8000  # I70: 5999.5 -> 8000.0
$ws.Cells.Item(70, 11).Value = 8000  # K70: 5999.5 -> 8000.0
$ws.Cells.Item(70, 13).Value = -7730  # M70: -5729.5 -> -7730.0
$ws.Cells.Item(73, 8).Value = 8000  # H73: 5999.5 -> 8000.0
$ws.Cells.Item(73, 9).Value = 8000  # I73: 5999.5 -> 8000.0
$ws.Cells.Item(73, 11).Value = 8000  # K73: 5999.5 -> 8000.0
$ws.Cells.Item(73, 13).Value = -7064  # M73: -5063.5 -> -7064.0
$ws.Cells.Item(122, 8).Value = 2864.6924  # H122: 3520.111 -> 2864.6924
$ws.Cells.Item(122, 9).Value = 2813.6667  # I122: 3393.1428 -> 2813.6667
$ws.Cells.Item(122, 10).Value = 2979.5  # J122: 3964.5 -> 2979.5
$ws.Cells.Item(122, 11).Value = 8441.000100000001  # K122: 10179.4284 -> 8441.000100000001
$ws.Cells.Item(122, 12).Value = 8938.5  # L122: 11893.5 -> 8938.5
$ws.Cells.Item(122, 13).Value = -5991.000100000001  # M122: -7729.428400000001 -> -5991.000100000001
$ws.Cells.Item(122, 14).Value = -13838.5  # N122: -16793.5 -> -13838.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 4610.4443  # H61: 5327.7144 -> 4610.4443
$ws.Cells.Item(61, 9).Value = 4532  # I61: 4131.6665 -> 4532.0
$ws.Cells.Item(61, 10).Value = 4649.6665  # J61: 6224.75 -> 4649.6665
$ws.Cells.Item(61, 11).Value = 4532  # K61: 4131.6665 -> 4532.0
$ws.Cells.Item(61, 12).Value = 4649.6665  # L61: 6224.75 -> 4649.6665
$ws.Cells.Item(61, 13).Value = -4330  # M61: -3929.6665 -> -4330.0
$ws.Cells.Item(61, 14).Value = -5053.6665  # N61: -6628.75 -> -5053.6665
$ws.Cells.Item(68, 8).Value = 0  # H68: 1333.0 -> 0.0
$ws.Cells.Item(68, 9).Value = 0  # I68: 1333.0 -> 0.0
$ws.Cells.Item(68, 11).Value = 0  # K68: 1333.0 -> 0.0
$ws.Cells.Item(68, 13).ClearContents()  # M68: -584.0 -> (removed)
$ws.Cells.Item(71, 8).Value = 0  # H71: 1333.0 -> 0.0
$ws.Cells.Item(71, 9).Value = 0  # I71: 1333.0 -> 0.0
$ws.Cells.Item(71, 11).Value = 0  # K71: 6665.0 -> 0.0
$ws.Cells.Item(71, 13).ClearContents()  # M71: -2921.0 -> (removed)
$ws.Cells.Item(113, 8).Value = 4610.4443  # H113: 5327.7144 -> 4610.4443
$ws.Cells.Item(113, 9).Value = 4532  # I113: 4131.6665 -> 4532.0
$ws.Cells.Item(113, 10).Value = 4649.6665  # J113: 6224.75 -> 4649.6665
$ws.Cells.Item(113, 11).Value = 4532  # K113: 4131.6665 -> 4532.0
$ws.Cells.Item(113, 12).Value = 4649.6665  # L113: 6224.75 -> 4649.6665
$ws.Cells.Item(113, 13).Value = -2362  # M113: -1961.6665 -> -2362.0
$ws.Cells.Item(113, 14).Value = -8989.6665  # N113: -10564.75 -> -8989.6665
$ws.Cells.Item(132, 8).Value = 4065  # H132: 4498.5557 -> 4065.0
$ws.Cells.Item(132, 9).Value = 4105.5713  # I132: 4498.5557 -> 4105.5713
$ws.Cells.Item(132, 10).Value = 3497  # J132: 0.0 -> 3497.0
$ws.Cells.Item(132, 11).Value = 12316.7139  # K132: 13495.6671 -> 12316.7139
$ws.Cells.Item(132, 12).Value = 10491  # L132: 0.0 -> 10491.0
$ws.Cells.Item(132, 13).Value = -9786.713899999999  # M132: -10965.6671 -> -9786.713899999999
$ws.Cells.Item(132, 14).Value = -15551  # N132: None -> -15551.0
$ws.Cells.Item(136, 8).Value = 2969.2856  # H136: 3130.8333 -> 2969.2856
$ws.Cells.Item(136, 9).Value = 2557  # I136: 2757.0 -> 2557.0
$ws.Cells.Item(136, 10).Value = 4000  # J136: 5000.0 -> 4000.0
$ws.Cells.Item(136, 11).Value = 7671  # K136: 8271.0 -> 7671.0
$ws.Cells.Item(136, 12).Value = 12000  # L136: 15000.0 -> 12000.0
$ws.Cells.Item(136, 13).Value = -5121  # M136: -5721.0 -> -5121.0
$ws.Cells.Item(136, 14).Value = -17100  # N136: -20100.0 -> -17100.0

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(13, 8).Value = 449  # H13: 449.5 -> 449.0
$ws.Cells.Item(13, 10).Value = 0  # J13: 450.0 -> 0.0
$ws.Cells.Item(13, 12).Value = 0  # L13: 450.0 -> 0.0
$ws.Cells.Item(13, 14).ClearContents()  # N13: -730.0 -> (removed)
$ws.Cells.Item(96, 8).Value = 998  # H96: 1165.3334 -> 998.0
$ws.Cells.Item(96, 10).Value = 900  # J96: 1200.0 -> 900.0
$ws.Cells.Item(96, 12).Value = 900  # L96: 1200.0 -> 900.0
$ws.Cells.Item(96, 14).Value = -3646  # N96: -3946.0 -> -3646.0
$ws.Cells.Item(126, 8).Value = 7998.222  # H126: 7999.0 -> 7998.222
$ws.Cells.Item(126, 10).Value = 7998.5713  # J126: 8000.0 -> 7998.5713
$ws.Cells.Item(126, 12).Value = 23995.7139  # L126: 24000.0 -> 23995.7139
$ws.Cells.Item(126, 14).Value = -28935.7139  # N126: -28940.0 -> -28935.7139
